# "made endpoint for intro text"
#
# The "Worksheet" sheet is a per-student reconciliation log: each student
# gets two rows (Trial 1 / Trial 2) recording whether they answered the
# intro-text question correctly, how long it took, and the date attempted.
#
# This edit:
#   - fixes the trial/date/elapsed-time bookkeeping for the existing
#     student (ikleiman@stonybrook.edu), whose two attempts had been
#     mis-recorded (duplicate Trial numbers, wrong elapsed time) and whose
#     date needs to move from 2019-10-21 to 2019-11-04
#   - appends three more students who have since completed the
#     "intro text" endpoint: chaotsai@stonybrook.edu (2019-12-03),
#     vlgarcia@stonybrook.edu (2019-12-03), pstdenis@stonybrook.edu
#     (2019-12-04)
#
# Values that could otherwise be auto-interpreted by Excel's input parser
# (the literal words "true"/"false" -> Boolean, "2019-11-04" -> a date
# serial) are entered with a leading apostrophe so they land as plain text,
# matching how this sheet already stores them everywhere else.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - ikleiman@stonybrook.edu, Trial 1, false, 2019-11-04
$ws.Range("D2").Value = 6
$ws.Range("E2").Value = "'2019-11-04"

# Row 3 - ikleiman@stonybrook.edu, Trial 1, false, 2019-11-04
$ws.Range("B3").Value = 1
$ws.Range("D3").Value = 6
$ws.Range("E3").Value = "'2019-11-04"

# Row 4 - ikleiman@stonybrook.edu, Trial 2, false, 2019-11-04
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = "'false"
$ws.Range("D4").Value = 13
$ws.Range("E4").Value = "'2019-11-04"

# Row 5 - chaotsai@stonybrook.edu, Trial 1, true, 2019-12-03
$ws.Range("A5").Value = "'chaotsai@stonybrook.edu"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "'true"
$ws.Range("D5").Value = 193
$ws.Range("E5").Value = "'2019-12-03"

# Row 6 - chaotsai@stonybrook.edu, Trial 2, false, 2019-12-03
$ws.Range("A6").Value = "'chaotsai@stonybrook.edu"
$ws.Range("B6").Value = 2
$ws.Range("D6").Value = 199
$ws.Range("E6").Value = "'2019-12-03"

# Row 7 - vlgarcia@stonybrook.edu, Trial 1, false, 2019-12-03
$ws.Range("A7").Value = "'vlgarcia@stonybrook.edu"
$ws.Range("B7").Value = 1
$ws.Range("D7").Value = 10
$ws.Range("E7").Value = "'2019-12-03"

# Row 8 (new) - vlgarcia@stonybrook.edu, Trial 2, false, 2019-12-03
$ws.Range("A8").Value = "'vlgarcia@stonybrook.edu"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "'false"
$ws.Range("D8").Value = 38
$ws.Range("E8").Value = "'2019-12-03"

# Row 9 (new) - pstdenis@stonybrook.edu, Trial 1, false, 2019-12-04
$ws.Range("A9").Value = "'pstdenis@stonybrook.edu"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "'false"
$ws.Range("D9").Value = 11
$ws.Range("E9").Value = "'2019-12-04"

# Row 10 (new) - pstdenis@stonybrook.edu, Trial 2, false, 2019-12-04
$ws.Range("A10").Value = "'pstdenis@stonybrook.edu"
$ws.Range("B10").Value = 2
$ws.Range("C10").Value = "'false"
$ws.Range("D10").Value = 19
$ws.Range("E10").Value = "'2019-12-04"

# The leading apostrophes above leave a "quote prefix" formatting hint on
# the cells; strip it back off so these rows keep the same (default, no
# explicit style) formatting as the rest of the sheet.
$ws.Range("A2:E10").ClearFormats()
